# "Generate Report for Archive"
#
# The localization-status report is regenerated: every cell that still read
# the old "Ready for handoff" status is now "In Translation", and the
# (now shorter) status column on each sheet is narrowed to match.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet: status is mirrored into the per-locale columns E/F ---
$wsOverview = $wb.Worksheets.Item("Overview")
$usedOverview = $wsOverview.UsedRange
for ($r = 1; $r -le $usedOverview.Rows.Count; $r++) {
    foreach ($c in 5, 6) {
        $cell = $wsOverview.Cells.Item($r, $c)
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Per-locale sheets: status lives in column C ---
foreach ($sheetName in "zh-cn", "de-de") {
    $ws = $wb.Worksheets.Item($sheetName)
    $used = $ws.UsedRange
    for ($r = 1; $r -le $used.Rows.Count; $r++) {
        $cell = $ws.Cells.Item($r, 3)
        if ($cell.Value() -eq $oldStatus) {
            $cell.Value = $newStatus
        }
    }
}

# --- Narrow the status columns now that the text is shorter ---
# The runtime quantizes ColumnWidth to a 1/6-character pixel grid, so 12.5
# characters is the input that lands closest to the regenerated report's
# target column width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wb.Worksheets.Item("zh-cn").Columns.Item(3).ColumnWidth = 12.5
$wb.Worksheets.Item("de-de").Columns.Item(3).ColumnWidth = 12.5
